$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right after "2021-Q4" (before "总计").
#    NOTE: Worksheets.Item(N) re-resolves by position every time it's
#    evaluated, so a handle captured *before* the insert would silently
#    start pointing at the new sheet afterwards. Fetch "总计" by name once
#    the sheet collection has its final shape instead of caching it early.
# ---------------------------------------------------------------------------
$sheetQ4 = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add($null, $sheetQ4)
$newSheet.Name = "2022-Q1"

$sheetTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2) Fill in the header row + 10 holding rows for "2022-Q1".
#    Values that look numeric but must be stored as TEXT (to match the
#    source sheet's "inlineStr" typing) are written with a leading
#    apostrophe so Excel keeps them as text instead of coercing to a number;
#    a later format-only paste strips the resulting quote-prefix style back
#    off again so the cell ends up with no explicit style, just like the
#    template sheet.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = [char]([int][char]'B' + $i)
    $newSheet.Range("$col`1").Value = $headers[$i]
}

$rows = @(
    @("011230", "创金合信数字经济主题股票C", "17.18", "92.17", "4.07", "0.6992", 5),
    @("011229", "创金合信数字经济主题股票A", "12.18", "92.17", "4.07", "0.4957", 5),
    @("000654", "华商新锐产业灵活配置混合", "13.02", "81.79", "3.24", "0.4218", 2),
    @("004423", "华商研究精选灵活配置混合", "9.97", "82.17", "3.24", "0.3230", 2),
    @("012491", "华商核心引力混合型证券投资基金A", "5.05", "83.13", "3.24", "0.1636", 2),
    @("008961", "华商科技创新混合", "2.87", "88.86", "3.25", "0.0933", 4),
    @("002504", "鹏华金鼎灵活配置混合A", "2.49", "77.53", "3.49", "0.0869", 10),
    @("012492", "华商核心引力混合型证券投资基金C", "0.51", "83.13", "3.24", "0.0165", 2),
    @("005161", "华商上游产业股票", "0.36", "89.02", "3.23", "0.0116", 7),
    @("002505", "鹏华金鼎灵活配置混合C", "0.26", "77.53", "3.49", "0.0091", 10)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $newSheet.Range("A$r").Value = $i
    $newSheet.Range("B$r").Value = "'" + $data[0]
    $newSheet.Range("C$r").Value = "'" + $data[1]
    $newSheet.Range("D$r").Value = "'" + $data[2]
    $newSheet.Range("E$r").Value = "'" + $data[3]
    $newSheet.Range("F$r").Value = "'" + $data[4]
    $newSheet.Range("G$r").Value = "'" + $data[5]
    $newSheet.Range("H$r").Value = $data[6]
}

# ---------------------------------------------------------------------------
# Normalize formatting to match the "2021-Q4" template sheet:
#  - row 1 (B1:H1) + column A index cells get the bold/bordered header style
#  - the rest of the data cells carry no explicit style
# This format-only paste also clears the quote-prefix style that got
# attached to the numeric-looking text cells above.
# ---------------------------------------------------------------------------
$sheetQ4.Range("A1:H2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$sheetQ4.Range("A2:H2").Copy()
$newSheet.Range("A3:H11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) "总计": insert a new row 2 ("2022-Q1" summary) above the existing
#    "2021-Q4" row (which shifts down to row 3).
# ---------------------------------------------------------------------------
$sheetTotal.Rows.Item(2).Insert()

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 10
$sheetTotal.Range("D2").Value = 2.32

# The "index" column (A) is a recomputed 0-based row counter, not a value
# that should merely shift down with the row -- bump the old row to 1.
$sheetTotal.Range("A3").Value = 1

# Normalize the new row's formatting off the still-correct row 3 below it.
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)
$sheetTotal.Range("B3:D3").Copy()
$sheetTotal.Range("B2:D2").PasteSpecial(-4122)

$excel.CutCopyMode = $false
